$wb = $excel.ActiveWorkbook

# --- Sheet: Intermediate for Mapping (columns V = Max Cr, W = Max Date) ---
$ws1 = $wb.Worksheets.Item("Intermediate for Mapping")

# These new values look like numbers/dates, so force the cells to stay
# text (matching the source workbook's inline-string / text storage)
# by applying a text number format before writing the value.
$sheet1TextCells = @("V4","W4","V8","W8","V13","W13","V14")
foreach ($addr in $sheet1TextCells) {
    $ws1.Range($addr).NumberFormat = "@"
}

$ws1.Range("V4").Value = "6.33"
$ws1.Range("W4").Value = "2001-09-07"

$ws1.Range("V8").Value = "5.32"
$ws1.Range("W8").Value = "2011-03-23"

$ws1.Range("V9").Value = "No Detect Data"
$ws1.Range("W9").Value = "No Detect Data"

$ws1.Range("V13").Value = "1.9"
$ws1.Range("W13").Value = "2005-03-24"

$ws1.Range("V14").Value = "6.01"

# --- Sheet: Intermediate Exhibit (columns G = Max Cr, H = Max Date) ---
$ws2 = $wb.Worksheets.Item("Intermediate Exhibit")

$sheet2TextCells = @("G6","H6","G11","H11","G17","H17","G18")
foreach ($addr in $sheet2TextCells) {
    $ws2.Range($addr).NumberFormat = "@"
}

$ws2.Range("G6").Value = "6.33"
$ws2.Range("H6").Value = "9/7/01"

$ws2.Range("G11").Value = "5.32"
$ws2.Range("H11").Value = "3/23/11"

$ws2.Range("G13").Value = "NA"
$ws2.Range("H13").Value = "NA"

$ws2.Range("G17").Value = "1.9"
$ws2.Range("H17").Value = "3/24/05"

$ws2.Range("G18").Value = "6.01"
